$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a cell's value as plain text, preserving the original
# "no explicit style" formatting (avoids Excel auto-converting
# number-looking strings like "7.03" or "1.20" into numeric values,
# and avoids leaving a stray text-format style behind).
function Set-CellText {
    param(
        [string]$CellRef,
        [string]$Text
    )
    $cell = $ws.Range($CellRef)
    $cell.NumberFormat = "@"
    $cell.Value = $Text
    $cell.Style = "Normal"
}

Set-CellText "D2" '41.224.60'
Set-CellText "E2" '  -0.27%  '
Set-CellText "D3" '2.187.31'
Set-CellText "E3" '  -1.47%  '
Set-CellText "E4" '  -0.16%  '
Set-CellText "D5" '256.04'
Set-CellText "E5" '  +4.89%  '
Set-CellText "E6" '  +0.00%  '
Set-CellText "D7" '68.53'
Set-CellText "E7" '  -1.62%  '
Set-CellText "E8" '  -0.13%  '
Set-CellText "D9" '0.573'
Set-CellText "E9" '  +2.83%  '
Set-CellText "D10" '37.20'
Set-CellText "E10" '  -5.07%  '
Set-CellText "D11" '58.85'
Set-CellText "E11" '  +1.16%  '
Set-CellText "D12" '0.0938'
Set-CellText "E12" '  -1.51%  '
Set-CellText "D13" '7.03'
Set-CellText "E13" '  +4.26%  '
Set-CellText "D14" '0.104'
Set-CellText "E14" '  +0.48%  '
Set-CellText "D15" '2.509.93'
Set-CellText "E15" '  -1.58%  '
Set-CellText "D16" '0.871'
Set-CellText "E16" '  +3.28%  '
Set-CellText "D17" '14.42'
Set-CellText "E17" '  -2.71%  '
Set-CellText "D18" '2.173.57'
Set-CellText "E18" '  -2.17%  '
Set-CellText "D19" '41.230.89'
Set-CellText "E19" '  -0.16%  '
Set-CellText "D20" '0.0₃0957'
Set-CellText "E20" '  +0.37%  '
Set-CellText "D21" '6.17'
Set-CellText "E21" '  +1.37%  '
Set-CellText "D22" '72.07'
Set-CellText "E22" '  -0.26%  '
Set-CellText "D23" '232.74'
Set-CellText "E23" '  +0.36%  '
Set-CellText "E24" '  -3.19%  '
Set-CellText "D25" '11.74'
Set-CellText "E25" '  +19.71%  '
Set-CellText "D26" '3.84'
Set-CellText "E26" '  +6.06%  '
Set-CellText "E27" '  -0.08%  '
Set-CellText "E28" '  +4.21%  '
Set-CellText "E29" '  -5.98%  '
Set-CellText "B30" 'Toncoin'
Set-CellText "C30" 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
Set-CellText "D30" '2.17'
Set-CellText "E30" '  -0.64%  '
Set-CellText "B31" 'Monero'
Set-CellText "C31" 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
Set-CellText "D31" '169.19'
Set-CellText "E31" '  -1.66%  '
Set-CellText "D32" '20.65'
Set-CellText "E32" '  +0.74%  '
Set-CellText "D34" '0.0749'
Set-CellText "E34" '  +4.67%  '
Set-CellText "E35" '  -0.47%  '
Set-CellText "D36" '5.48'
Set-CellText "E36" '  +4.52%  '
Set-CellText "D37" '26.41'
Set-CellText "E37" '  +9.14%  '
Set-CellText "D38" '4.17'
Set-CellText "E38" '  +6.82%  '
Set-CellText "D39" '4.59'
Set-CellText "E39" '  -0.28%  '
Set-CellText "D40" '0.0298'
Set-CellText "E40" '  +7.55%  '
Set-CellText "D41" '2.20'
Set-CellText "E41" '  -3.50%  '
Set-CellText "D42" '12.25'
Set-CellText "E42" '  +13.64%  '
Set-CellText "D43" '5.66'
Set-CellText "E43" '  -3.19%  '
Set-CellText "D44" '63.58'
Set-CellText "D45" '4.98'
Set-CellText "E45" '  -0.15%  '
Set-CellText "D46" '0.198'
Set-CellText "E46" '  -3.22%  '
Set-CellText "D47" '8.61'
Set-CellText "E47" '  -2.59%  '
Set-CellText "E48" '  +0.62%  '
Set-CellText "D49" '1.20'
Set-CellText "E49" '  +9.33%  '
Set-CellText "E50" '  +0.21%  '
Set-CellText "E51" '  +0.23%  '

